$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "29.395.05"
Set-TextValue "E2" "  +0.09%  "
Set-TextValue "D3" "1.841.92"
Set-TextValue "E3" "  -0.22%  "
Set-TextValue "D4" "0.9990"
Set-TextValue "E4" "  +0.08%  "
Set-TextValue "D5" "239.40"
Set-TextValue "E5" "  -0.41%  "
Set-TextValue "D6" "0.6275"
Set-TextValue "E6" "  +0.02%  "
Set-TextValue "E7" "  +0.09%  "
Set-TextValue "E8" "  -0.67%  "
Set-TextValue "D9" "0.2896"
Set-TextValue "E9" "  -0.17%  "
Set-TextValue "D10" "24.89"
Set-TextValue "E10" "  +1.69%  "
Set-TextValue "D11" "0.07709"
Set-TextValue "E11" "  -0.31%  "
Set-TextValue "D12" "1.843.36"
Set-TextValue "E12" "  -0.15%  "
Set-TextValue "D13" "4.973"
Set-TextValue "E13" "  -0.49%  "
Set-TextValue "D14" "0.6762"
Set-TextValue "E14" "  -0.77%  "
Set-TextValue "E15" "  -3.01%  "
Set-TextValue "D16" "81.85"
Set-TextValue "E16" "  -0.31%  "
Set-TextValue "D17" "6.245"
Set-TextValue "E17" "  +0.91%  "
Set-TextValue "D18" "29.425.95"
Set-TextValue "E18" "  +0.05%  "
Set-TextValue "D19" "233.35"
Set-TextValue "E19" "  +1.74%  "
Set-TextValue "E20" "  -0.02%  "
Set-TextValue "E21" "  +0.13%  "
Set-TextValue "D22" "7.324"
Set-TextValue "E22" "  -2.09%  "
Set-TextValue "E23" "  +0.18%  "
Set-TextValue "D24" "158.27"
Set-TextValue "E24" "  -0.81%  "
Set-TextValue "D25" "8.495"
Set-TextValue "E25" "  +0.79%  "
Set-TextValue "D26" "0.1351"
Set-TextValue "E26" "  -1.69%  "
Set-TextValue "D27" "17.36"
Set-TextValue "E27" "  -0.98%  "
Set-TextValue "D28" "0.07093"
Set-TextValue "E28" "  +7.99%  "
Set-TextValue "D29" "1.463"
Set-TextValue "E29" "  +2.80%  "
Set-TextValue "E30" "  -0.17%  "
Set-TextValue "B31" "Filecoin"
Set-TextValue "C31" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D31" "4.040"
Set-TextValue "E31" "  -1.36%  "
Set-TextValue "B32" "InternetComputer(DFINITY)"
Set-TextValue "C32" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D32" "4.047"
Set-TextValue "E32" "  -1.25%  "
Set-TextValue "D33" "1.820"
Set-TextValue "E33" "  -0.70%  "
Set-TextValue "D34" "1.140"
Set-TextValue "E34" "  -0.18%  "
Set-TextValue "D35" "0.7006"
Set-TextValue "E35" "  +0.44%  "
Set-TextValue "D36" "2.575"
Set-TextValue "E36" "  -0.17%  "
Set-TextValue "D37" "6.996"
Set-TextValue "E37" "  +2.83%  "
Set-TextValue "D38" "0.01837"
Set-TextValue "E38" "  +0.34%  "
Set-TextValue "D39" "2.814"
Set-TextValue "E39" "  -0.88%  "
Set-TextValue "D40" "1.238.20"
Set-TextValue "E40" "  -2.21%  "
Set-TextValue "D41" "0.9584"
Set-TextValue "E41" "  +5.32%  "
Set-TextValue "D42" "1.000"
Set-TextValue "E42" "  +0.14%  "
Set-TextValue "D43" "1.995.75"
Set-TextValue "E43" "  -0.59%  "
Set-TextValue "D44" "101.03"
Set-TextValue "E44" "  -0.20%  "
Set-TextValue "D45" "65.47"
Set-TextValue "E45" "  -1.32%  "
Set-TextValue "E46" "  +0.74%  "
Set-TextValue "D47" "1.729"
Set-TextValue "E47" "  -0.89%  "
Set-TextValue "D48" "6.968"
Set-TextValue "E48" "  -1.56%  "
Set-TextValue "D49" "8.946"
Set-TextValue "E49" "  -0.95%  "
Set-TextValue "E50" "  -2.36%  "
Set-TextValue "D51" "0.3900"
Set-TextValue "E51" "  -1.51%  "
